$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
$ws.Range("A1").Value = "Email Id"
$ws.Range("B1").Value = "Password"

# --- Update credential rows with refreshed test data (all four accounts now
#     share the one "Selenium@123" password) ---
$ws.Range("A2").Value = "rafselenium1@gmail.com"
$ws.Range("B2").Value = "Selenium@123"

$ws.Range("A3").Value = "rafselenium2@gmail.com"
$ws.Range("B3").Value = "Selenium@123"

$ws.Range("A4").Value = "rafselenium3@yahoo.com"
$ws.Range("B4").Value = "Selenium@123"

$ws.Range("A5").Value = "rafselenium4@yahoo.com"
$ws.Range("B5").Value = "Selenium@123"

# Row 6 is emptied out
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()

# --- Rebuild hyperlinks: relationship ids/targets are unchanged, but which
#     cell each one is attached to changes (B6's old hyperlink now lives on
#     A4 instead) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:rafiasultana12345@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Selenium@123")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:rafiasultana122@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Selenium@123")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:sharmin@123")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:rashidmohammed@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:mohammed@123")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:rashida@123")

# Re-applying a hyperlink nudges formatting, so restore the table's original
# "hyperlink + thin border" look uniformly across A2:B5 (this now also
# covers A4, which previously used a plain bordered style).
$ws.Range("B3").Copy()
$ws.Range("A2:B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 6 keeps the hyperlink font but loses its border
$ws.Range("A6:B6").Style = "Hyperlink"
$ws.Range("A6:B6").Borders.LineStyle = -4142

# --- New empty rows below the table (hyperlink-styled, no border) ---
$ws.Range("A9:B13").Style = "Hyperlink"

# Select A9 as the new active cell (matches the saved selection state)
$ws.Range("A9").Select()
